# Update registration/visitor counts in column F across the sheets of
# "广州-漫展信息.xlsx" to reflect the latest scraped totals.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 279    # was 277
$ws.Range("F4").Value = 973    # was 970
$ws.Range("F5").Value = 247    # was 246
$ws.Range("F7").Value = 671    # was 667
$ws.Range("F10").Value = 10    # was 7
$ws.Range("F12").Value = 182   # was 180
$ws.Range("F13").Value = 43    # was 41
$ws.Range("F14").Value = 771   # was 768
$ws.Range("F16").Value = 1912  # was 1908
$ws.Range("F17").Value = 425   # was 423
$ws.Range("F18").Value = 5915  # was 5858
$ws.Range("F19").Value = 422   # was 420
$ws.Range("F21").Value = 36    # was 34
$ws.Range("F23").Value = 8     # was 7

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 501    # was 499

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 362    # was 361
$ws.Range("F4").Value = 353    # was 348

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 362    # was 361
$ws.Range("F6").Value = 353    # was 348
$ws.Range("F7").Value = 279    # was 277
$ws.Range("F12").Value = 501   # was 499
$ws.Range("F13").Value = 501   # was 499
$ws.Range("F14").Value = 973   # was 970
$ws.Range("F17").Value = 247   # was 246
$ws.Range("F19").Value = 671   # was 667
$ws.Range("F23").Value = 10    # was 7
$ws.Range("F25").Value = 182   # was 180
$ws.Range("F27").Value = 43    # was 41
$ws.Range("F29").Value = 771   # was 768
$ws.Range("F32").Value = 1912  # was 1908
$ws.Range("F33").Value = 425   # was 423
$ws.Range("F34").Value = 5915  # was 5858
$ws.Range("F36").Value = 422   # was 420
$ws.Range("F38").Value = 36    # was 34
$ws.Range("F41").Value = 8     # was 7
